$wb = $excel.ActiveWorkbook

# Update "covariate_importance" sheet (sheet 1)
$ws1 = $wb.Worksheets.Item("covariate_importance")

$covariateData = @(
    @("region", 100, 100, 100),
    @("state", 100, 100, 100),
    @("percenttwoormoreraces", 95, 86, 86),
    @("percentwhite", 88, 64, 64),
    @("percentasian", 57.99999999999999, 31, 31),
    @("percentfreelunchqualified", 61, 29, 29),
    @("percentblackorafricanamerican", 59, 22, 22),
    @("rplthemes", 39, 13, 13),
    @("percenthispaniclatino", 50, 11, 11),
    @("percentamericanindianoralaskanative", 38, 10, 10),
    @("schoollevel", 43, 8, 8),
    @("derivedtotalenrolled", 33, 7.000000000000001, 7.000000000000001),
    @("cntycaseschange", 45, 6, 6),
    @("locale", 30, 6, 6),
    @("percentnativehawaiianorotherpacificislander", 47, 6, 6),
    @("percentnotspecified", 0, 1, 0)
)

$row = 2
foreach ($entry in $covariateData) {
    $ws1.Cells.Item($row, 1).Value = $entry[0]
    $ws1.Cells.Item($row, 2).Value = $entry[1]
    $ws1.Cells.Item($row, 3).Value = $entry[2]
    $ws1.Cells.Item($row, 4).Value = $entry[3]
    $row++
}

# Update "strategy_importance" sheet (sheet 2)
$ws2 = $wb.Worksheets.Item("strategy_importance")

$strategyData = @(
    @("contacttracing", 99, 99, 99),
    @("hvacsystems", 94, 93, 93),
    @("cleaning", 89, 85, 89),
    @("screeningtestingforstudents", 67, 64, 67),
    @("masks", 44, 38, 42),
    @("physicaldistancing", 16, 9, 13),
    @("vaccination", 14, 4, 9),
    @("hepafilters", 5, 4, 3),
    @("quarantine", 3, 3, 3),
    @("stayhome", 1, 1, 1)
)

$row = 2
foreach ($entry in $strategyData) {
    $ws2.Cells.Item($row, 1).Value = $entry[0]
    $ws2.Cells.Item($row, 2).Value = $entry[1]
    $ws2.Cells.Item($row, 3).Value = $entry[2]
    $ws2.Cells.Item($row, 4).Value = $entry[3]
    $row++
}
